$d = $word.ActiveDocument

# Locate the three paragraphs to remove: the blank paragraph, the
# "Ver no Jupiter..." paragraph, and the "(c) 2020..." paragraph that sit
# between the page-break paragraph following "Requisitos" and the final
# blank/page-break paragraphs at the end of the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ver no Jupiter*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startPara = $target.Previous()   # the blank paragraph right before "Ver no Jupiter..."
    $endPara = $target.Next()         # the "(c) 2020 ..." paragraph right after it

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}

Write-Output "done"
